# Refresh the cryptos price/volume table (rows 2-51) with the latest scrape.
# A handful of Price cells (column D) are plain-looking decimals (e.g. "0.9969",
# "0.05340") that Excel would otherwise auto-coerce to numbers on assignment
# (losing trailing zeros / switching to scientific notation). For those we
# force the cell to Text via NumberFormat "@" before writing the value, then
# restore the cell style to "Normal" so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.328.40'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '1.712.68'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9969'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9978'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.22%  '
$ws.Range("E7").Value = '  -1.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2582'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06168'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.10%  '
$ws.Range("D10").Value = '1.712.25'
$ws.Range("E10").Value = '  -1.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06940'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.61%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.46'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.60%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.462'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -3.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5956'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.64%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '76.32'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9974'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9969'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.31%  '
$ws.Range("D18").Value = '26.242.55'
$ws.Range("E18").Value = '  -1.11%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007088'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.29%  '
$ws.Range("E20").Value = '  -2.90%  '
$ws.Range("D21").Value = '1.934.22'
$ws.Range("E21").Value = '  -1.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.395'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.397'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.029'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.29%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '136.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.92%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.13'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.393'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.70%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.724'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.39%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '105.21'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.47%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.863'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.63%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07929'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.54%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.599'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.31%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04425'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.65%  '
$ws.Range("B34").Value = 'Frax'
$ws.Range("C34").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9966'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("B35").Value = 'HuobiToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.599'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.40%  '
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9873'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.21%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6161'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.41%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9321'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.04%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.984'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.72%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.367'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("B41").Value = 'PaxDollar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9966'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.61%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01468'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.67%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '99.83'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.54%  '
$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.382'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.45%  '
$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.3797'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.83%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '6.842'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.79%  '
$ws.Range("B47").Value = 'Algorand'
$ws.Range("C47").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1147'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.35%  '
$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.05340'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '30.50'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.32%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.672'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '50.97'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.79%  '
